# Skill.xlsx edit: rename RequireDistance/TargetType columns to
# Melee/EffectObjType, change the "Melee" column's data type from
# float distance (2.5) to an int flag (0), and add hp for heroes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): column names -----------------------------------
$ws.Range("Q1").Value = "Melee"
$ws.Range("S1").Value = "EffectObjType"

# --- Row 2: field "Type" metadata for the Melee column (float -> int) ---
$ws.Range("Q2").Value = "int"

# --- Row 10: Chinese field descriptions ----------------------------------
$ws.Range("Q10").Value = "Melee:0;Ranger:1"
$ws.Range("S10").Value = "作用效果目标`n0敌人`n1自己`n2队友(include self)"

# --- Data rows 11-61: Melee column values 2.5 -> 0, re-highlighted ------
# Column P (CoolDownTime) row 11 already carries the yellow-fill /
# thin-border style that should now cover the whole Melee column, so
# copy its format down over Q11:Q61 (mirrors a Format Painter pass)
# before overwriting the values.
$ws.Range("P11").Copy()
$ws.Range("Q11:Q61").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 11; $r -le 61; $r++) {
    $ws.Cells.Item($r, 17).Value = 0
}

# --- Restore the on-screen selection to where the edit left off ---------
$ws.Range("S2").Select()
